# Updated cryptos list on Mon Apr 22 14:47:51 UTC 2024 with GitHub Actions
# Applies the latest price/volume refresh to the cryptos worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.279.13"
$ws.Range("D3").Value = "3.197.72"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.11%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.198.24"
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.514"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.52%  "
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.17%  "
$ws.Range("D15").Value = "3.721.54"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "66.229.30"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("E17").Value = "  +4.93%  "
$ws.Range("D18").Value = "3.201.23"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.111"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.738"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "  +2.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.97%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.92%  "
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.99%  "
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0901"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "485.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.64%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("E42").Value = "  +5.03%  "
$ws.Range("E43").Value = "  +6.35%  "
$ws.Range("E44").Value = "  +12.34%  "
$ws.Range("D45").Value = "2.923.01"
$ws.Range("E45").Value = "  -4.16%  "
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.13%  "
